$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
